$d = $word.ActiveDocument

# 1. Title heading and matching bold text near the end (same text, both occurrences
#    are replaced in one pass since Replace is wdReplaceAll (2) over the whole $d.Content range)
$d.Content.Find.Execute("Play Gold Bonanza for Free: Review and Game details", $true, $false, $false, $false, $false, $true, 1, $false, "Play Gold Bonanza | Free Slot Game", 2) | Out-Null

# 2. "What we like" bullet list rotation
$d.Content.Find.Execute("Special collection function offers added rewards", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting special features like Shifting Wilds and Free Spins", 2) | Out-Null
$d.Content.Find.Execute("Interesting paytable and payouts", $true, $false, $false, $false, $false, $true, 1, $false, "Well-crafted graphics and interesting gold rush theme", 2) | Out-Null
$d.Content.Find.Execute("Well-crafted graphics and theme", $true, $false, $false, $false, $false, $true, 1, $false, "Special collection function for additional rewards", 2) | Out-Null

# 3. "What we don't like" bullet list wording tweaks
$d.Content.Find.Execute("RTP is not the best at 95.45%", $true, $false, $false, $false, $false, $true, 1, $false, "Not the best RTP at 95.45%", 2) | Out-Null
$d.Content.Find.Execute("The gold rush theme may not be for everyone", $true, $false, $false, $false, $false, $true, 1, $false, "Theme may not be appealing to everyone", 2) | Out-Null

# 4. Meta description (italic) paragraph near the end
$d.Content.Find.Execute("Discover the special collection function and high volatility of Gold Bonanza, a well-designed slot game with a good payout potential. Play for free now.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Gold Bonanza, an exciting online slot game to play for free.", 2) | Out-Null
